# Work Profile and new tenant support
# Appends new sprint/interview-history data rows to the AMSIN, BETA and AMS
# sheets, and refreshes the (previously slightly stale) last row of AMS.

function Add-SprintRow {
    param($ws, $row, $dateText, $timeSerial, $sprintName, $total, $pass, $fail, $timeTaken, $formatSourceRow)

    # Columns C-G: plain values. When written into a brand new row these
    # naturally pick up the worksheet's column default style (matches the
    # neighbouring rows).
    $ws.Cells.Item($row, 3).Value = $sprintName
    $ws.Cells.Item($row, 4).Value = $total
    $ws.Cells.Item($row, 5).Value = $pass
    $ws.Cells.Item($row, 6).Value = $fail
    $ws.Cells.Item($row, 7).Value = $timeTaken

    # Column A holds a date-like label that must stay literal text (not get
    # silently converted into a date serial number) - prefix with a quote to
    # force text entry, then copy the neighbouring cell's formatting across
    # so the quote-prefix styling doesn't linger.
    $ws.Cells.Item($row, 1).Value = "'" + $dateText
    $ws.Cells.Item($formatSourceRow, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    # Column B holds the run's date-time serial number and must keep the
    # "YYYY-MM-DD HH:MM:SS" number format used throughout the column.
    $ws.Cells.Item($row, 2).Value = $timeSerial
    $ws.Cells.Item($formatSourceRow, 2).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# AMSIN sheet: add rows 71-74
# ---------------------------------------------------------------------
$wsAMSIN = $wb.Worksheets.Item("AMSIN")
Add-SprintRow $wsAMSIN 71 "2023-03-10" 44995.80265944445 "174ffiinnalrun" 165 161 4 4.93 70
Add-SprintRow $wsAMSIN 72 "2023-03-13" 44998.47240579861 "174finalrun" 165 161 4 4.17 71
Add-SprintRow $wsAMSIN 73 "2023-03-31" 45016.47215266203 "175fnlrun" 165 156 9 7.96 72
Add-SprintRow $wsAMSIN 74 "2023-04-12" 45028.7219215049 "176fstrtail" 165 165 0 4.3 73

# ---------------------------------------------------------------------
# BETA sheet: add rows 33-34
# ---------------------------------------------------------------------
$wsBETA = $wb.Worksheets.Item("BETA")
Add-SprintRow $wsBETA 33 "2023-03-13" 44998.54521554398 "174beta" 165 165 0 4.79 32
Add-SprintRow $wsBETA 34 "2023-03-31" 45016.5515381713 "175beta" 165 165 0 4.77 33

# ---------------------------------------------------------------------
# AMS sheet: refresh row 36 (tiny timestamp correction + style refresh)
# and add rows 37-40
# ---------------------------------------------------------------------
$wsAMS = $wb.Worksheets.Item("AMS")
$wsAMS.Range("A36:G36").Clear()
Add-SprintRow $wsAMS 36 "2023-02-20" 44977.84599061342 "live173" 165 165 0 5.45 35
Add-SprintRow $wsAMS 37 "2023-03-01" 44986.70179755787 "173angularvrs" 165 165 0 4.81 36
Add-SprintRow $wsAMS 38 "2023-03-02" 44987.44341578703 "liveangular173" 165 165 0 4.85 37
Add-SprintRow $wsAMS 39 "2023-03-13" 44998.84790481481 "174live" 165 165 0 4.89 38
Add-SprintRow $wsAMS 40 "2023-03-31" 45016.80041184028 "175live" 165 165 0 4.82 39

Write-Host "Applied Work Profile and new tenant support updates"
